# Weekly update: insert a new weekly price record for Plátano
# (Agrícola del Norte S.A. de Arica) ahead of the existing history,
# shifting the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 189; Excel shifts rows 189:197 down to 190:198
# and copies formatting from the row above (keeps the date style on column D).
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row 189 with the new weekly record.
$ws.Cells.Item(189, 1).Value  = 1
$ws.Cells.Item(189, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(189, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(189, 4).Value  = 44610
$ws.Cells.Item(189, 5).Value  = 15
$ws.Cells.Item(189, 6).Value  = "Fruta"
$ws.Cells.Item(189, 7).Value  = 100108
$ws.Cells.Item(189, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(189, 9).Value  = 100108006
$ws.Cells.Item(189, 10).Value = "Plátano"
$ws.Cells.Item(189, 11).Value = "Sin especificar"
$ws.Cells.Item(189, 12).Value = "Pintón"
$ws.Cells.Item(189, 13).Value = 130
$ws.Cells.Item(189, 14).Value = 17000
$ws.Cells.Item(189, 15).Value = 18000
$ws.Cells.Item(189, 16).Value = 17500
$ws.Cells.Item(189, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(189, 18).Value = "Ecuador"
$ws.Cells.Item(189, 19).Value = 875
$ws.Cells.Item(189, 20).Value = 20
